$d = $word.ActiveDocument

# --- Content edit ------------------------------------------------------
# Update the memory-size sentence in one shot:
#   "...64 GB of system memory, ... 512 GB or more."
#   -> "...128 GB of system memory, ... 512 GB or even more."
# (Kept as a single Find/Replace call on this paragraph -- issuing a second,
# separate Find/Replace over the same paragraph would cause the engine to
# re-flow/merge all of that paragraph's runs, wiping out pre-existing run
# boundaries that must stay put, e.g. around "Third".)
$old = "64 GB of system memory, while a professional desktop can have 512 GB or more."
$new = "128 GB of system memory, while a professional desktop can have 512 GB or even more."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
    $true, 1, $false, $new, 2)

# --- Recreate the run boundaries / "_GoBack" bookmark that a live, manual
#     edit in Word leaves behind. Zero-length bookmarks split the enclosing
#     run without altering any text, and (unlike Find/Replace) do not
#     trigger a paragraph-wide run re-flow. Do this after the text edit.

# Isolate "128" into its own run (mirrors typing "128" over the old "64").
$r128 = $d.Content
$r128.Find.Execute("128", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$s128 = $r128.Start
$e128 = $r128.End
$d.Bookmarks.Add("TempA", $d.Range($s128, $s128))
$d.Bookmarks("TempA").Delete()
$d.Bookmarks.Add("TempB", $d.Range($e128, $e128))
$d.Bookmarks("TempB").Delete()

# Isolate "even " into its own run (mirrors typing "even " right before
# "more.").
$rEven = $d.Content
$rEven.Find.Execute("even ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$sEven = $rEven.Start
$eEven = $rEven.End
$d.Bookmarks.Add("TempD", $d.Range($sEven, $sEven))
$d.Bookmarks("TempD").Delete()

# Isolate "more." from the following " Last, ..." text.
$rMore = $d.Content
$rMore.Find.Execute("more.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$sMore = $rMore.Start
$eMore = $rMore.End
$d.Bookmarks.Add("TempC", $d.Range($eMore, $eMore))
$d.Bookmarks("TempC").Delete()

# Word automatically maintains a hidden "_GoBack" bookmark at the location of
# the most recent edit. Re-create it (as a zero-length bookmark) right
# before "more." (i.e. right after "even ", where text was last inserted),
# replacing wherever it previously sat in the document.
$d.Bookmarks.Add("_GoBack", $d.Range($sMore, $sMore))
